# Revert "Powerpoint writer: consolidate text run nodes."
#
# The consolidated writer had merged "word" + "trailing space" into a single
# run (e.g. "A " or "Just "). This restores the original shape where each
# word and each inter-word space live in their own <a:r> run, while leaving
# the run's (empty) formatting (<a:rPr/>) and the overall visible text
# untouched.
#
# Mechanism: re-assigning (even unchanged) text to a `Characters(start,
# length)` sub-range of a TextRange makes PowerPoint's text engine carve that
# sub-range out into its own run, splitting whatever run it used to share
# text with. Doing this once per inter-word space is exactly the "un-merge"
# we need.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

function Split-TrailingSpaces($textRange, [int[]]$spacePositions) {
    foreach ($pos in $spacePositions) {
        $textRange.Characters($pos, 1).Text = " "
    }
}

# --- Title placeholder: "A " + "slide" -> "A" + " " + "slide" ---
$title = Get-ShapeByName $s "Title 1"
$titleRange = $title.TextFrame.TextRange
Split-TrailingSpaces $titleRange @(2)

# --- Caption textbox: "Just " + "an " + "image " + "on " + "this " + "side"
#     -> "Just" + " " + "an" + " " + "image" + " " + "on" + " " + "this" +
#        " " + "side"
# Text is "Just an image on this side"; the spaces that need to become their
# own runs sit right after "Just", "an", "image", "on" and "this".
$caption = Get-ShapeByName $s "TextBox 3"
$captionRange = $caption.TextFrame.TextRange
Split-TrailingSpaces $captionRange @(5, 8, 14, 17, 22)
